# Add a new "2022-Q4" sheet right after "总计" and before "2021-Q4",
# populate it with the new quarter's holdings data, and update the
# "总计" (summary) sheet with a new leading row for 2022-Q4 (shifting
# the existing rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet before "2021-Q4".
# ---------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Re-fetch fresh references (avoid any stale COM handles after Add()).
$q4_2021 = $wb.Worksheets.Item("2021-Q4")
$q4_2022 = $wb.Worksheets.Item("2022-Q4")

# Copy the header/format layout from the existing "2021-Q4" sheet so the
# new sheet matches the established look (bold/centered header row,
# styled index column, etc.). Column A of row 1 is intentionally left
# untouched/empty, matching every other quarter sheet in this workbook.
$q4_2021.Range("B1:H1").Copy() | Out-Null
$q4_2022.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$q4_2021.Range("A2").Copy() | Out-Null
$q4_2022.Range("A2").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------
# 2. Fill in the "2022-Q4" sheet contents.
# ---------------------------------------------------------------
$q4_2022.Range("B1").Value = "基金代码"
$q4_2022.Range("C1").Value = "基金名称"
$q4_2022.Range("D1").Value = "基金规模"
$q4_2022.Range("E1").Value = "股票总仓位"
$q4_2022.Range("F1").Value = "仓位占比"
$q4_2022.Range("G1").Value = "持有市值(亿元)"
$q4_2022.Range("H1").Value = "仓位排名"

$q4_2022.Range("A2").Value = 0

# Numeric-looking values that must stay TEXT (not be coerced to numbers).
# Force text via NumberFormat "@" while assigning, then drop the format
# override again afterwards so the cell ends up unstyled, matching the
# rest of the data rows (only the index column A carries a style).
$q4_2022.Range("B2").NumberFormat = "@"
$q4_2022.Range("B2").Value = "501089"
$q4_2022.Range("B2").ClearFormats()

$q4_2022.Range("C2").Value = "方正富邦消费红利指数增强（LOF）"

$q4_2022.Range("D2").NumberFormat = "@"
$q4_2022.Range("D2").Value = "0.22"
$q4_2022.Range("D2").ClearFormats()

$q4_2022.Range("E2").NumberFormat = "@"
$q4_2022.Range("E2").Value = "94.07"
$q4_2022.Range("E2").ClearFormats()

$q4_2022.Range("F2").NumberFormat = "@"
$q4_2022.Range("F2").Value = "3.99"
$q4_2022.Range("F2").ClearFormats()

$q4_2022.Range("G2").NumberFormat = "@"
$q4_2022.Range("G2").Value = "0.0088"
$q4_2022.Range("G2").ClearFormats()

$q4_2022.Range("H2").Value = 4

# ---------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: shift rows 2-6 down to 3-7,
#    then write the new 2022-Q4 figures into row 2.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Give the new row 7 the same formatting as row 6 before writing into it.
$summary.Range("A6").Copy() | Out-Null
$summary.Range("A7").PasteSpecial(-4122) | Out-Null

for ($r = 7; $r -ge 3; $r--) {
    $src = $r - 1
    $summary.Range("B$r").Value = $summary.Range("B$src").Value2
    $summary.Range("C$r").Value = $summary.Range("C$src").Value2
    $summary.Range("D$r").Value = $summary.Range("D$src").Value2
}

# Re-sequence the index column (A2:A7 = 0..5).
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01
